$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest observation (row 2): A, B, C, D columns shift up by one row,
# matching the naive forecaster's corrected window alignment.
$ws.Rows("2:2").Delete()

# Column E ("y_1_forecast") is recomputed for the corrected window. The first
# four remaining rows no longer carry a forecast value; clear any values that
# shifted into E2:E5 as a result of the row deletion above.
$ws.Range("E2:E5").ClearContents()

# Recalculated y_1_forecast values for rows 6-18.
$ws.Range("E6").Value = 0.9288717675470126
$ws.Range("E7").Value = 1.416624765035412
$ws.Range("E8").Value = 1.075154359849861
$ws.Range("E9").Value = 1.214249019249602
$ws.Range("E10").Value = 1.180565832117297
$ws.Range("E11").Value = 1.638669199130427
$ws.Range("E12").Value = 1.540918326052476
$ws.Range("E13").Value = 0.5773070399857971
$ws.Range("E14").Value = -0.415982961498651
$ws.Range("E15").Value = 1.426719405738508
$ws.Range("E16").Value = 0.5205511175203181
$ws.Range("E17").Value = 0.3517304536567734
$ws.Range("E18").Value = 0.4186921370205043
